$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.947.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.643.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.36%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5090"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.28%  "

$ws.Range("E7").Value = "  +0.24%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2570"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.23%  "

$ws.Range("E9").Value = "  +0.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.66%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07794"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.11%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.310"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.95%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.648.26"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5462"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅7860"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.74%  "

$ws.Range("E16").Value = "  +1.13%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.008.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.46%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.005"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.20%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "197.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.12%  "

$ws.Range("E20").Value = "  +1.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.965"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.044"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.008"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.871"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.38%  "

$ws.Range("E26").Value = "  +0.63%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.898"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.66%  "

$ws.Range("E28").Value = "  +0.70%  "

$ws.Range("E29").Value = "  -0.39%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05033"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.264"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.196"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.02%  "

$ws.Range("E33").Value = "  +0.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.365"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.35%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.8954"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.47%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.593"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.37%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.133.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.48%  "

$ws.Range("E38").Value = "  -1.38%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01558"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.39%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₈132"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +14.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.006"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.546"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.51%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.644"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.72%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8175"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.38%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "100.03"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.30%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.779.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4541"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.35%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.007"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.34%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05080"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.008"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.45%  "

